$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.51"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.98%"
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.91"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.54%"
$ws.Range("E3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.423"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.58%"
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08090"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.56%"
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.723"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.98%"
$ws.Range("E6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.902"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.30%"
$ws.Range("E7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.301"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.96%"
$ws.Range("E8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.751"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.28%"
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9412"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.82%"
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1212"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.40%"
$ws.Range("E11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1901"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.87%"
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09523"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.63%"
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04167"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "8.39%"
$ws.Range("E14").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.78%"
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001275"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.48%"
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006098"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.16%"
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.571"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.22%"
$ws.Range("E18").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.493"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.10%"
$ws.Range("E20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1356"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.45%"
$ws.Range("E21").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.77%"
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04379"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.70%"
$ws.Range("E23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001239"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.48%"
$ws.Range("E24").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004316"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.42%"
$ws.Range("E25").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.51%"
$ws.Range("E26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004011"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.49%"
$ws.Range("E27").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02677"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.76%"
$ws.Range("E39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05440"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.56%"
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007804"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.07%"
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009761"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.43%"
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1393"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.91%"
$ws.Range("E43").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.65%"
$ws.Range("E44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009961"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-14.72%"
$ws.Range("E45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007069"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.99%"
$ws.Range("E46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.49%"
$ws.Range("E47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003468"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.08%"
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002282"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.15%"
$ws.Range("E49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002111"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.49%"
$ws.Range("E50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002011"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.49%"
$ws.Range("E51").ClearFormats()
